$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document starts as 5 one-line paragraphs: "1" (+ the hidden _GoBack
# bookmark), "2", "3", "4", "" (empty).  The target shape is 8 paragraphs:
#   1: "R" / "htxsjfuyjtjdufyjdyt"              (was "1", keeps its paraId)
#   2: "F" / "jhty" (+ _GoBack bookmark)          (new paragraph)
#   3: ""                                          (was "2", text removed, keeps paraId)
#   4: "1"                                         (new paragraph)
#   5: "2"                                         (new paragraph)
#   6: "3"                                         (untouched)
#   7: "4"                                         (untouched)
#   8: "" -> "5"                                   (keeps paraId, run added)
#
# Plain InsertAfter()/Text= typing gets auto-coalesced into a single <w:r>
# whenever the formatting matches, which would not reproduce the separate
# runs the diff shows. Range.InsertXML, however, splices raw OOXML in
# verbatim (no run-merging), so each paragraph below is rebuilt by replacing
# its whole Range (text + trailing paragraph mark) with literal OOXML.
#
# Edits are applied from the end of the story backwards so that earlier
# offsets stay valid while later ones move.
# ---------------------------------------------------------------------------

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$w14Ns = "xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

# --- Paragraph 5 (last, empty -> "5") --------------------------------------
$p5 = $d.Paragraphs.Item(5)
$xml5 = "<w:p $wNs $w14Ns w14:paraId='1EADC3CA' w14:textId='77777777' w:rsidR='00353123' w:rsidRPr='000A15A1' w:rsidRDefault='00353123'>" +
        "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>5</w:t></w:r>" +
        "</w:p>"
$p5.Range.InsertXML($xml5)

# --- Paragraph 2 ("2" -> empty, then two new paragraphs "1" and "2") -------
$p2 = $d.Paragraphs.Item(2)
$xml2 = "<w:p $wNs $w14Ns w14:paraId='088ADA5A' w14:textId='77777777' w:rsidR='00353123' w:rsidRDefault='00353123'>" +
        "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
        "</w:p>" +
        "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>1</w:t></w:r>" +
        "</w:p>" +
        "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>2</w:t></w:r>" +
        "</w:p>"
$p2.Range.InsertXML($xml2)

# --- Paragraph 1 ("1" + bookmark -> "R"/"htxsjfuyjtjdufyjdyt", then a new --
# --- paragraph "F"/"jhty" that inherits the bookmark) ----------------------
$p1 = $d.Paragraphs.Item(1)
$xml1 = "<w:p $wNs $w14Ns w14:paraId='72F3A842' w14:textId='77777777' w:rsidR='00F40B3D' w:rsidRDefault='00353123'>" +
        "<w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>R</w:t></w:r>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>htxsjfuyjtjdufyjdyt</w:t></w:r>" +
        "</w:p>" +
        "<w:p $wNs><w:pPr><w:rPr><w:lang w:val='en-US'/></w:rPr></w:pPr>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>F</w:t></w:r>" +
        "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>jhty</w:t></w:r>" +
        "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
        "</w:p>"
$p1.Range.InsertXML($xml1)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host $i ":" $d.Paragraphs.Item($i).Range.Text
}
